$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Questions"
$ws.Range("B1").Value = "Answers"
$ws.Range("A2").Value = "Who is the CEO of Intel?"
$ws.Range("B2").Value = "Patrick P. Gelsinger"

$ws.Range("B6").Select()
